# Update countries & provincias Spain
# Applies the 25-Jun-2020 15:28 data refresh to the "Pais" ranking sheet.
#
# Most rows keep their country label and simply get refreshed totals, but a
# few countries crossed over each other in the case-count ranking, so their
# row positions (and therefore the country label shown on that row) swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp footer (row 1) -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 15:28"

# --- Rows whose country stays put, only the numbers refresh -------------
$ws.Range("B4").Value  = 2463923
$ws.Range("C4").Value  = 1369
$ws.Range("E4").Value  = 1299007
$ws.Range("G4").Value  = 27
$ws.Range("H4").Value  = 124308

$ws.Range("B7").Value  = 474587
$ws.Range("C7").Value  = 1602
$ws.Range("D7").Value  = 272636
$ws.Range("E7").Value  = 187036

$ws.Range("B18").Value = 170639
$ws.Range("C18").Value = 3372
$ws.Range("D18").Value = 117882
$ws.Range("E18").Value = 51329
$ws.Range("G18").Value = 41
$ws.Range("H18").Value = 1428

$ws.Range("D27").Value = 16890
$ws.Range("E27").Value = 34391

$ws.Range("B62").Value = 13372
$ws.Range("C62").Value = 137
$ws.Range("D62").Value = 12154
$ws.Range("E62").Value = 954
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 264

$ws.Range("B63").Value = 12636
$ws.Range("C63").Value = 21
$ws.Range("D63").Value = 11460
$ws.Range("E63").Value = 573

$ws.Range("B80").Value = 5595
$ws.Range("C80").Value = 150
$ws.Range("D80").Value = 2166
$ws.Range("E80").Value = 3164
$ws.Range("G80").Value = 6
$ws.Range("H80").Value = 265

$ws.Range("B93").Value = 3796
$ws.Range("C93").Value = 120
$ws.Range("D93").Value = 2322
$ws.Range("E93").Value = 1299
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 175

$ws.Range("B129").Value = 1056
$ws.Range("C129").Value = 5
$ws.Range("D129").Value = 917
$ws.Range("E129").Value = 72

$ws.Range("D132").Value = 562
$ws.Range("E132").Value = 429

$ws.Range("B142").Value = 788
$ws.Range("C142").Value = 26
$ws.Range("D142").Value = 221
$ws.Range("E142").Value = 562

$ws.Range("B162").Value = 272
$ws.Range("C162").Value = 7
$ws.Range("D162").Value = 161
$ws.Range("E162").Value = 104

# --- Ranking overtakes: Uzbekistan's two neighbours refresh in place ----
$ws.Range("B76").Value = 7087
$ws.Range("C76").Value = 186
$ws.Range("D76").Value = 4724
$ws.Range("E76").Value = 2343
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 20

# Paises Bajos overtakes Argentina for rank 36 -> row 32 now shows
# "Paises Bajos" with refreshed numbers; Argentina drops to row 33 with its
# previous (unchanged) totals.
$ws.Range("A32").Value = "Paises Bajos"
$ws.Range("B32").Value = 49914
$ws.Range("C32").Value = 110
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 6100

$ws.Range("A33").Value = "Argentina"
$ws.Range("B33").Value = 49851
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 13816
$ws.Range("E33").Value = 34919
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 1116

# Consejo Danes para los Refugiados overtakes Senegal for rank 81 -> row 77
# now shows "Consejo Danes para los Refugiados" with refreshed numbers;
# Senegal drops to row 78 with its previous (unchanged) totals.
$ws.Range("A77").Value = "Consejo Danes para los Refugiados"
$ws.Range("B77").Value = 6411
$ws.Range("C77").Value = 198
$ws.Range("D77").Value = 885
$ws.Range("E77").Value = 5384
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 142

$ws.Range("A78").Value = "Senegal"
$ws.Range("B78").Value = 6233
$ws.Range("C78").Value = 104
$ws.Range("D78").Value = 4162
$ws.Range("E78").Value = 1977
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 94

# --- Tied-count rows that simply swapped tie-break order (no value change) ---
# Fiyi now ranks just ahead of Dominica (both on 18 total cases).
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# Groenlandia now ranks just ahead of Islas Malvinas (both on 13 total cases).
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"
